$d = $word.ActiveDocument

$replacements = @(
    @("91×61=5551", "77×24=1848"),
    @("21×52=1092", "58×14=812"),
    @("26×12=312",  "64×43=2752"),
    @("34×27=918",  "90×76=6840"),
    @("54×25=1350", "78×42=3276"),
    @("78×12=936",  "13×66=858"),
    @("75×76=5700", "51×58=2958"),
    @("72×80=5760", "77×82=6314"),
    @("35×25=875",  "47×66=3102"),
    @("31×28=868",  "18×15=270"),
    @("77×29=2233", "97×51=4947"),
    @("75×46=3450", "70×45=3150"),
    @("87×78=6786", "97×53=5141"),
    @("51×60=3060", "42×85=3570"),
    @("56×21=1176", "51×18=918"),
    @("49×23=1127", "65×78=5070"),
    @("28×43=1204", "96×73=7008"),
    @("89×12=1068", "17×85=1445"),
    @("32×56=1792", "51×24=1224"),
    @("81×82=6642", "28×53=1484"),
    @("31×77=2387", "12×95=1140"),
    @("82×44=3608", "81×68=5508"),
    @("73×96=7008", "48×48=2304"),
    @("58×55=3190", "24×63=1512"),
    @("67×15=1005", "69×21=1449")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
